$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new rows before row 2 (shifts existing rows 2-13 down to 11-22)
$ws.Range("A2:A10").EntireRow.Insert()

# The newly inserted rows inherited formatting from the row above; restore the
# plain (unstyled) A/B columns and the date-formatted C/D columns by copying
# the formats from the (now-shifted) row that used to be row 3.
$ws.Range("A12:D12").Copy()
$ws.Range("A2:D10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# New period rows to place in rows 2-10, 'orden' 1-9
$newData = @(
    @(201303, 1, 41456, 41639),
    @(201401, 2, 41640, 41820),
    @(201403, 3, 41821, 42004),
    @(201501, 4, 42005, 42185),
    @(201503, 5, 42186, 42369),
    @(201601, 6, 42370, 42551),
    @(201603, 7, 42552, 42735),
    @(201701, 8, 42736, 42916),
    @(201703, 9, 42917, 43100)
)

for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = 2 + $i
    $row = $newData[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Renumber the 'orden' column (B) for the rows that were shifted down (now rows 11-22)
for ($r = 11; $r -le 22; $r++) {
    $ws.Cells.Item($r, 2).Value = $r - 1
}

# Update selection to match the final state
$ws.Range("J13").Select()
